$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the 3 remaining data rows (7, 8, 9) with the new shortage items ---
# Row 7: ORGASOL LIGHT CREAM
$ws.Range("C7").Value = "ORGASOL LIGHT CREAM"
$ws.Range("H7").Value = "-1:0"
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = "130.00"
$ws.Range("P7").Value = "130.0000"
$ws.Range("Q7").Value = "1:0"

# Row 8: PRISBRINA  CAPS
$ws.Range("C8").Value = "PRISBRINA  CAPS"
$ws.Range("H8").Value = "-1:0"
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = "150.00"
$ws.Range("P8").Value = "150.0000"
$ws.Range("Q8").Value = "1:0"

# Row 9: QUICK NAIL  LOTION
$ws.Range("C9").Value = "QUICK NAIL  LOTION"
$ws.Range("H9").Value = "-23:0"
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = "85.00"
$ws.Range("P9").Value = "1955.0000"
$ws.Range("Q9").Value = "23:0"

# --- Remove the 5 rows that no longer have products (old rows 10-14) ---
$ws.Rows("10:14").Delete()

# --- Update the new totals row (was row 15, now row 10) ---
$ws.Range("P10").Value = 2235

# --- Update the timestamp in the footer (was row 16, now row 11) ---
$ws.Range("A11").Value = "Saturday, 24 May, 2025 10:04 AM"
